# Update numeric counts (column F) across the four worksheets
# ("展览", "演出", "本地生活", "全部类型") to reflect the latest
# scraped data (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1295
$ws.Range("F4").Value = 146
$ws.Range("F5").Value = 889
$ws.Range("F6").Value = 24
$ws.Range("F7").Value = 1044
$ws.Range("F8").Value = 148
$ws.Range("F9").Value = 1079
$ws.Range("F10").Value = 821
$ws.Range("F11").Value = 40
$ws.Range("F12").Value = 699
$ws.Range("F13").Value = 1360
$ws.Range("F14").Value = 1031
$ws.Range("F16").Value = 748
$ws.Range("F21").Value = 644
$ws.Range("F22").Value = 1245
$ws.Range("F23").Value = 162
$ws.Range("F25").Value = 858
$ws.Range("F26").Value = 5218
$ws.Range("F27").Value = 263
$ws.Range("F30").Value = 5837
$ws.Range("F32").Value = 992
$ws.Range("F34").Value = 62
$ws.Range("F47").Value = 352

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F13").Value = 128
$ws.Range("F17").Value = 784
$ws.Range("F18").Value = 16
$ws.Range("F32").Value = 149
$ws.Range("F33").Value = 104
$ws.Range("F36").Value = 67
$ws.Range("F41").Value = 490

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 754
$ws.Range("F6").Value = 381

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1295
$ws.Range("F6").Value = 381
$ws.Range("F10").Value = 146
$ws.Range("F11").Value = 889
$ws.Range("F13").Value = 24
$ws.Range("F14").Value = 1044
$ws.Range("F15").Value = 148
$ws.Range("F16").Value = 1079
$ws.Range("F17").Value = 821
$ws.Range("F18").Value = 40
$ws.Range("F20").Value = 699
$ws.Range("F21").Value = 1360
$ws.Range("F23").Value = 1031
$ws.Range("F24").Value = 748
$ws.Range("F29").Value = 644
$ws.Range("F30").Value = 1245
$ws.Range("F31").Value = 162
$ws.Range("F34").Value = 5218
$ws.Range("F35").Value = 263
$ws.Range("F38").Value = 5837
$ws.Range("F40").Value = 992
$ws.Range("F43").Value = 62
$ws.Range("F47").Value = 67
$ws.Range("F50").Value = 490
